# Auto-generated Excel COM-interop script to apply scheduled-runner market data refresh
# to the per-job (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) leve-profit tables.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4488.846
$ws.Range("I76").Value = 4395.7144
$ws.Range("K76").Value = 4395.7144
$ws.Range("M76").Value = -4080.7144
$ws.Range("H79").Value = 4488.846
$ws.Range("I79").Value = 4395.7144
$ws.Range("K79").Value = 4395.7144
$ws.Range("M79").Value = -3303.7144
$ws.Range("H86").Value = 32132.908
$ws.Range("I86").Value = 43510.707
$ws.Range("J86").Value = 1792.1111
$ws.Range("K86").Value = 43510.707
$ws.Range("L86").Value = 1792.1111
$ws.Range("M86").Value = -42387.707
$ws.Range("N86").Value = -4038.1111
$ws.Range("H89").Value = 32132.908
$ws.Range("I89").Value = 43510.707
$ws.Range("J89").Value = 1792.1111
$ws.Range("K89").Value = 217553.535
$ws.Range("L89").Value = 8960.5555
$ws.Range("M89").Value = -211937.535
$ws.Range("N89").Value = -20192.5555
$ws.Range("H111").Value = 7267
$ws.Range("I111").Value = 8528.352999999999
$ws.Range("K111").Value = 25585.059
$ws.Range("M111").Value = -22518.059
$ws.Range("H116").Value = 17666.5
$ws.Range("I116").Value = 17666.5
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 17666.5
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -14224.5
$ws.Range("N116").ClearContents()
$ws.Range("H137").Value = 1007.625
$ws.Range("I137").Value = 966.4865
$ws.Range("J137").Value = 1146
$ws.Range("K137").Value = 2899.4595
$ws.Range("L137").Value = 3438
$ws.Range("M137").Value = -349.4594999999999
$ws.Range("N137").Value = -8538
$ws.Range("H138").Value = 2192.6736
$ws.Range("I138").Value = 2052.2856
$ws.Range("J138").Value = 2248.8286
$ws.Range("K138").Value = 6156.8568
$ws.Range("L138").Value = 6746.485799999999
$ws.Range("M138").Value = -1016.8568
$ws.Range("N138").Value = -17026.4858

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2158.7334
$ws.Range("I63").Value = 1819.4445
$ws.Range("K63").Value = 1819.4445
$ws.Range("M63").Value = -1133.4445
$ws.Range("H66").Value = 2158.7334
$ws.Range("I66").Value = 1819.4445
$ws.Range("K66").Value = 9097.2225
$ws.Range("M66").Value = -5665.2225
$ws.Range("H74").Value = 1271.9642
$ws.Range("I74").Value = 1371.3158
$ws.Range("J74").Value = 1062.2222
$ws.Range("K74").Value = 1371.3158
$ws.Range("L74").Value = 1062.2222
$ws.Range("M74").Value = -497.3158000000001
$ws.Range("N74").Value = -2810.2222
$ws.Range("H77").Value = 1271.9642
$ws.Range("I77").Value = 1371.3158
$ws.Range("J77").Value = 1062.2222
$ws.Range("K77").Value = 6856.579000000001
$ws.Range("L77").Value = 5311.111
$ws.Range("M77").Value = -2488.579000000001
$ws.Range("N77").Value = -14047.111
$ws.Range("H88").Value = 1691.2
$ws.Range("I88").Value = 1502
$ws.Range("J88").Value = 1975
$ws.Range("K88").Value = 1502
$ws.Range("L88").Value = 1975
$ws.Range("M88").Value = -1096
$ws.Range("N88").Value = -2787
$ws.Range("H91").Value = 1691.2
$ws.Range("I91").Value = 1502
$ws.Range("J91").Value = 1975
$ws.Range("K91").Value = 1502
$ws.Range("L91").Value = 1975
$ws.Range("M91").Value = -98
$ws.Range("N91").Value = -4783
$ws.Range("H102").Value = 64801.125
$ws.Range("I102").Value = 113607.555
$ws.Range("K102").Value = 113607.555
$ws.Range("M102").Value = -111985.555
$ws.Range("H122").Value = 2090.077
$ws.Range("I122").Value = 2132
$ws.Range("J122").Value = 2054.1428
$ws.Range("K122").Value = 6396
$ws.Range("L122").Value = 6162.428400000001
$ws.Range("M122").Value = -3946
$ws.Range("N122").Value = -11062.4284
$ws.Range("H132").Value = 4274.4165
$ws.Range("J132").Value = 4272.8887
$ws.Range("L132").Value = 12818.6661
$ws.Range("N132").Value = -17878.6661

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1352.7
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1352.7
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1352.7
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2254.7
$ws.Range("H132").Value = 2821.2222
$ws.Range("I132").Value = 2841.1785
$ws.Range("K132").Value = 8523.5355
$ws.Range("M132").Value = -5993.5355

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2107.4443
$ws.Range("I68").Value = 1099.6666
$ws.Range("J68").Value = 2611.3333
$ws.Range("K68").Value = 3298.9998
$ws.Range("L68").Value = 7833.999899999999
$ws.Range("M68").Value = -2487.9998
$ws.Range("N68").Value = -9455.999899999999
$ws.Range("H71").Value = 2107.4443
$ws.Range("I71").Value = 1099.6666
$ws.Range("J71").Value = 2611.3333
$ws.Range("K71").Value = 9896.999400000001
$ws.Range("L71").Value = 23501.9997
$ws.Range("M71").Value = -5840.999400000001
$ws.Range("N71").Value = -31613.9997
$ws.Range("H131").Value = 825.66
$ws.Range("J131").Value = 834.71136
$ws.Range("L131").Value = 2504.13408
$ws.Range("N131").Value = -12584.13408
$ws.Range("H132").Value = 1865.5834
$ws.Range("I132").Value = 1051.4
$ws.Range("J132").Value = 2447.1428
$ws.Range("K132").Value = 9462.6
$ws.Range("L132").Value = 22024.2852
$ws.Range("M132").Value = -6932.6
$ws.Range("N132").Value = -27084.2852
$ws.Range("H140").Value = 6986.4443
$ws.Range("I140").Value = 6986.4443
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 20959.3329
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -15779.3329
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 21306
$ws.Range("I141").Value = 25382.5
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 76147.5
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -70967.5
$ws.Range("N141").Value = -25360

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 160718.16
$ws.Range("I70").Value = 255464.5
$ws.Range("J70").Value = 9124
$ws.Range("K70").Value = 255464.5
$ws.Range("L70").Value = 9124
$ws.Range("M70").Value = -255194.5
$ws.Range("N70").Value = -9664
$ws.Range("H73").Value = 160718.16
$ws.Range("I73").Value = 255464.5
$ws.Range("J73").Value = 9124
$ws.Range("K73").Value = 255464.5
$ws.Range("L73").Value = 9124
$ws.Range("M73").Value = -254528.5
$ws.Range("N73").Value = -10996
$ws.Range("H80").Value = 91004480
$ws.Range("I80").Value = 200207360
$ws.Range("K80").Value = 200207360
$ws.Range("M80").Value = -200206362
$ws.Range("H83").Value = 91004480
$ws.Range("I83").Value = 200207360
$ws.Range("K83").Value = 1001036800
$ws.Range("M83").Value = -1001031808
$ws.Range("H126").Value = 1930
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1930
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 5790
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -10730

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6468.5
$ws.Range("I7").Value = 4250
$ws.Range("J7").Value = 7023.125
$ws.Range("K7").Value = 4250
$ws.Range("L7").Value = 7023.125
$ws.Range("M7").Value = -4138
$ws.Range("N7").Value = -7247.125
$ws.Range("H40").Value = 66043.25
$ws.Range("I40").Value = 147914.58
$ws.Range("J40").Value = 2365.5557
$ws.Range("K40").Value = 147914.58
$ws.Range("L40").Value = 2365.5557
$ws.Range("M40").Value = -147778.58
$ws.Range("N40").Value = -2637.5557
$ws.Range("H46").Value = 422229.6
$ws.Range("I46").Value = 332.5
$ws.Range("J46").Value = 844126.7
$ws.Range("K46").Value = 332.5
$ws.Range("L46").Value = 844126.7
$ws.Range("M46").Value = -144.5
$ws.Range("N46").Value = -844502.7
$ws.Range("H61").Value = 1981.1177
$ws.Range("I61").Value = 2025.625
$ws.Range("J61").Value = 1941.5555
$ws.Range("K61").Value = 2025.625
$ws.Range("L61").Value = 1941.5555
$ws.Range("M61").Value = -1823.625
$ws.Range("N61").Value = -2345.5555
$ws.Range("H113").Value = 1981.1177
$ws.Range("I113").Value = 2025.625
$ws.Range("J113").Value = 1941.5555
$ws.Range("K113").Value = 2025.625
$ws.Range("L113").Value = 1941.5555
$ws.Range("M113").Value = 144.375
$ws.Range("N113").Value = -6281.5555
$ws.Range("H126").Value = 6468.5
$ws.Range("I126").Value = 4250
$ws.Range("J126").Value = 7023.125
$ws.Range("K126").Value = 12750
$ws.Range("L126").Value = 21069.375
$ws.Range("M126").Value = -10280
$ws.Range("N126").Value = -26009.375

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 28593
$ws.Range("J129").Value = 28593
$ws.Range("L129").Value = 28593
$ws.Range("N129").Value = -38593
